# Auto update stock data
# Updates the date column (A) from 2025/12/17 -> 2025/12/18 and refreshes
# several metric values for specific rows in the active worksheet.
# Leading apostrophes force Excel to keep these as text values (matching
# the source data, which stores dates/numbers as text strings) instead of
# auto-converting them to a date serial number / numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ A = "'2025/12/18"; B = "'5.75" }
    8  = @{ A = "'2025/12/18"; B = "'8.18" }
    14 = @{ A = "'2025/12/18"; B = "'3.01" }
    20 = @{ A = "'2025/12/18"; B = "'13.04" }
    26 = @{ A = "'2025/12/18"; B = "'10.90" }
    32 = @{ A = "'2025/12/18"; B = "'27.30" }
    38 = @{ A = "'2025/12/18" }
    44 = @{ A = "'2025/12/18"; B = "'10.54" }
    50 = @{ A = "'2025/12/18"; B = "'11.42" }
    56 = @{ A = "'2025/12/18"; B = "'28.93" }
    62 = @{ A = "'2025/12/18"; B = "'10.84"; C = "'2.50"; D = "'6.30"; E = "'0.99" }
    68 = @{ A = "'2025/12/18"; B = "'13.08" }
    74 = @{ A = "'2025/12/18"; B = "'15.59" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
